# Fix Excel emoji problem: replace the "statut" (column A) marker values.
# Old -> New mapping (per shared-string table order):
#   📘 -> ⚠️
#   📕 -> -3
#   📙 -> +3
#   📗 -> ✅

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📙" = "+3"
    "📗" = "✅"
}

# "-3"/"+3" look like numbers, so Excel would silently coerce them into
# numeric cells on assignment. Prefix those replacements with an apostrophe
# so Excel keeps (and stores) them as plain text, exactly like the target
# shared-string entries.
$textForce = @{
    "-3" = $true
    "+3" = $true
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1
$firstRow = $usedRange.Row

# Column A holds the "statut" values starting on row 2 (row 1 is the header).
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $newVal = $map[$val]
        if ($textForce.ContainsKey($newVal)) {
            $cell.Value2 = "'" + $newVal
        } else {
            $cell.Value2 = $newVal
        }
    }
}
